# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# Map of row -> new value on the "展览" sheet (column F)
$exhibitUpdates = @{
    3  = 367
    4  = 168
    5  = 210
    6  = 359
    8  = 2213
    9  = 376
    10 = 5428
    11 = 124
    12 = 359
}

foreach ($row in $exhibitUpdates.Keys) {
    $sheetExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Map of row -> new value on the "全部类型" sheet (column F)
$allTypesUpdates = @{
    4  = 367
    5  = 168
    6  = 210
    7  = 359
    11 = 2213
    12 = 376
    13 = 5428
    14 = 124
    15 = 359
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
